# Update crypto price (D) and 1h volume/change (E) columns with refreshed values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.503.57"
$ws.Range("E2").Value = "  +1.75%  "
$ws.Range("D3").Value = "3.095.62"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'527.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.20%  "
$ws.Range("D6").Value = "'142.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.05%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +1.33%  "
$ws.Range("D9").Value = "'7.33"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("E10").Value = "  +0.16%  "
$ws.Range("E11").Value = "  +2.52%  "
$ws.Range("D12").Value = "3.628.83"
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").Value = "'26.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.08%  "
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").Value = "58.529.92"
$ws.Range("E16").Value = "  +1.65%  "
$ws.Range("D17").Value = "3.096.27"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").Value = "'6.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").Value = "'12.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.08%  "
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("D21").Value = "'341.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.94%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("E23").Value = "  +0.57%  "
$ws.Range("D24").Value = "'66.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").Value = "0.0₃0914"
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("D28").Value = "'6.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.00%  "
$ws.Range("D29").Value = "'7.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.89%  "
$ws.Range("E30").Value = "  +2.88%  "
$ws.Range("D31").Value = "'1.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.35%  "
$ws.Range("D32").Value = "'20.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("D33").Value = "'154.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("D34").Value = "'4.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.37%  "
$ws.Range("D35").Value = "'6.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.90%  "
$ws.Range("D36").Value = "'26.87"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.01%  "
$ws.Range("E37").Value = "  +3.54%  "
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("D39").Value = "3.139.49"
$ws.Range("E39").Value = "  +0.36%  "
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("D42").Value = "'36.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("E43").Value = "  +7.88%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "2.296.93"
$ws.Range("E45").Value = "  +0.40%  "
$ws.Range("D46").Value = "'0.0257"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("D47").Value = "'20.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.85%  "
$ws.Range("D48").Value = "'0.964"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.53%  "
$ws.Range("D49").Value = "'5.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.64%  "
$ws.Range("D50").Value = "'268.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.37%  "
$ws.Range("D51").Value = "'0.745"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.47%  "
